$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Min Score"
$ws.Range("C2").Value = "Max Score"

$ws.Range("G4").Select()
